$d = $word.ActiveDocument

$old = " à l’aide du site : http://www.passwordmeter.com"
$new = "."

$range = $d.Content
$found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

Write-Host "Found and replaced: $found"
